$d = $word.ActiveDocument

# The GPA listed in the cover letter was updated from 3.71 to 3.73.
$d.Content.Find.Execute("GPA of 3.71", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GPA of 3.73", 2) | Out-Null

# Word keeps the "_GoBack" bookmark (the last-edit marker) at the location
# of the most recent edit. Since the GPA digit was the last thing typed,
# move the zero-length "_GoBack" bookmark to sit right after the new value.
$find = $d.Content.Find
$find.Text = "GPA of 3.73"
$find.Execute() | Out-Null
$afterGpa = $find.Parent
$afterGpa.Collapse(0)
$d.Bookmarks.Add("_GoBack", $afterGpa) | Out-Null
